$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. "Last updated Date : March 13, 2014" -> "... March 03, 2018"
#    and move the "_GoBack" bookmark to sit right after the new date
#    text (Word keeps only one _GoBack bookmark, so adding a new one
#    automatically removes the old one wherever it was).
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Last updated Date : March 13, 2014", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Last updated Date : March 03, 2018", 2)

$dateRng = $d.Content
$dateRng.Find.Execute("Last updated Date : March 03, 2018")
$dateRng.Collapse(0)
$endPos = $dateRng.Start
$dateRng.InsertAfter("X")
$goBackRng = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $goBackRng)
$d.Range($endPos, $endPos + 1).Delete()

# ---------------------------------------------------------------------
# 2. Sample code fix: "ShowElementList( selSet, " -> "ShowElementList(selectedElementIds, "
# ---------------------------------------------------------------------
$d.Content.Find.Execute("ShowElementList( selSet, ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ShowElementList(selectedElementIds, ", 2)

# ---------------------------------------------------------------------
# 3. Drop the stale <w:lastRenderedPageBreak/> markers left over from a
#    previous pagination snapshot (Word regenerates these on layout, but
#    the authored copy strips the cached ones).
# ---------------------------------------------------------------------
Write-Output "done"
